$d = $word.ActiveDocument

# Update the title/date line
$d.Content.Find.Execute("2024-11-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-14 Thursday", 2) | Out-Null

# Update each table cell value (processed in row-major document order to avoid ambiguity from duplicate old values)
$newValues = @(
    "77-52=",
    "7+40=",
    "12+48=",
    "88-79=",
    "70+10=",
    "81-51=",
    "0+70=",
    "23+18=",
    "40+8=",
    "55-35=",
    "30+0=",
    "71-49=",
    "95-69=",
    "55-22=",
    "60-2=",
    "22+22=",
    "70-25=",
    "74-29=",
    "51+19=",
    "32-19=",
    "28-21=",
    "3+50=",
    "28+58=",
    "21+19=",
    "8+68=",
    "74-55=",
    "0+80=",
    "40+52=",
    "40+10=",
    "47-31=",
    "64-59=",
    "31+56=",
    "46-15=",
    "43-14=",
    "55-9=",
    "35+18=",
    "34+0=",
    "43-27=",
    "37+26=",
    "12+7=",
    "92-0=",
    "4+47=",
    "7+48=",
    "30-0=",
    "95-15=",
    "4+74=",
    "52+33=",
    "98-27=",
    "84+11=",
    "75-4=",
    "96-8=",
    "77-74=",
    "86-76=",
    "39+6=",
    "60+1=",
    "24+72=",
    "72-41=",
    "32+61=",
    "20+45=",
    "48+27=",
    "28+22=",
    "71+22=",
    "74-20=",
    "39+35=",
    "43-5=",
    "96-42=",
    "32+41=",
    "58-44=",
    "77-7=",
    "35+6=",
    "85-20=",
    "32+6=",
    "39-25=",
    "34-5=",
    "99-5=",
    "41+44=",
    "51-49=",
    "91-18=",
    "28-5=",
    "14+37=",
    "61+20=",
    "30+62=",
    "90-72=",
    "70-55=",
    "22+32=",
    "28+49=",
    "76+19=",
    "44-42=",
    "7+45=",
    "95-88=",
    "20-13=",
    "33+42=",
    "59+14=",
    "56-9=",
    "99-85=",
    "60-2=",
    "88-74=",
    "40+37=",
    "87-60=",
    "13-10="
)

$t = $d.Tables.Item(1)
$numCols = $t.Columns.Count
$idx = 0
foreach ($row in 1..$t.Rows.Count) {
    foreach ($col in 1..$numCols) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated " + $idx + " cells")